$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 12 entry (first so its string lands first in the shared string table)
$ws.Range("B12").Value = "쭈꾸미불고기"
$ws.Range("D12").Value = 45000
$ws.Range("E12").Value = "V"

# New row 11 entry
$ws.Range("B11").Value = "동연"
$ws.Range("C11").Value = 50000

# Rename contributors
$ws.Range("B6").Value = "영서"
$ws.Range("B4").Value = "형주"

# C17 gets a couple of spaces
$ws.Range("C17").Value = "  "

# Update the active selection to K5
$ws.Range("K5").Select()

$wb.Save()
